$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.00581
$ws.Range("H2").Value = 3.01743
$ws.Range("I2").Value = 0.003799625168827527
$ws.Range("J2").Value = 0.003799625168827527
$ws.Range("O2").Value = 0.3717075934090293
$ws.Range("P2").Value = 0.3717075934090293
$ws.Range("Q2").Value = 0.05059794259
$ws.Range("R2").Value = 0.45538148331
$ws.Range("S2").Value = 0.001412349527361257
$ws.Range("T2").Value = 0.001412349527361257

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.00581
$ws.Range("H3").Value = 3.01743
$ws.Range("I3").Value = 0.003799625168827527
$ws.Range("J3").Value = 0.003799625168827527
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.08503100000000001
$ws.Range("N3").Value = 0.255093
$ws.Range("O3").Value = 0.6282924065909707
$ws.Range("P3").Value = 0.6282924065909707
$ws.Range("Q3").Value = 0.08552503011000001
$ws.Range("R3").Value = 0.7697252709900001
$ws.Range("S3").Value = 0.002387275641466271
$ws.Range("T3").Value = 0.002387275641466271

# Row 4
$ws.Range("I4").Value = 0.9594121222074437
$ws.Range("J4").Value = 0.9594121222074438
$ws.Range("O4").Value = 0.3717075934090293
$ws.Range("P4").Value = 0.3717075934090293
$ws.Range("S4").Value = 0.3566207710331784
$ws.Range("T4").Value = 0.3566207710331785

# Row 5
$ws.Range("I5").Value = 0.9594121222074437
$ws.Range("J5").Value = 0.9594121222074438
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.08503100000000001
$ws.Range("N5").Value = 0.255093
$ws.Range("O5").Value = 0.6282924065909707
$ws.Range("P5").Value = 0.6282924065909707
$ws.Range("Q5").Value = 21.59522242163967
$ws.Range("R5").Value = 194.357001794757
$ws.Range("S5").Value = 0.6027913511742653
$ws.Range("T5").Value = 0.6027913511742654

# Row 6
$ws.Range("G6").Value = 9.336668333333334
$ws.Range("H6").Value = 28.010005
$ws.Range("I6").Value = 0.03527091597053946
$ws.Range("J6").Value = 0.03527091597053946
$ws.Range("O6").Value = 0.3717075934090293
$ws.Range("P6").Value = 0.3717075934090293
$ws.Range("Q6").Value = 0.4696873249538889
$ws.Range("R6").Value = 4.227185924585
$ws.Range("S6").Value = 0.01311046729274132
$ws.Range("T6").Value = 0.01311046729274132

# Row 7
$ws.Range("G7").Value = 9.336668333333334
$ws.Range("H7").Value = 28.010005
$ws.Range("I7").Value = 0.03527091597053946
$ws.Range("J7").Value = 0.03527091597053946
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.08503100000000001
$ws.Range("N7").Value = 0.255093
$ws.Range("O7").Value = 0.6282924065909707
$ws.Range("P7").Value = 0.6282924065909707
$ws.Range("Q7").Value = 0.7939062450516668
$ws.Range("R7").Value = 7.145156205465001
$ws.Range("S7").Value = 0.02216044867779814
$ws.Range("T7").Value = 0.02216044867779814

# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.4016586666666667
$ws.Range("H8").Value = 1.204976
$ws.Range("I8").Value = 0.001517336653189343
$ws.Range("J8").Value = 0.001517336653189343
$ws.Range("O8").Value = 0.3717075934090293
$ws.Range("P8").Value = 0.3717075934090293
$ws.Range("Q8").Value = 0.02020570699911111
$ws.Range("R8").Value = 0.181851362992
$ws.Range("S8").Value = 0.0005640055557483215
$ws.Range("T8").Value = 0.0005640055557483215

# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.4016586666666667
$ws.Range("H9").Value = 1.204976
$ws.Range("I9").Value = 0.001517336653189343
$ws.Range("J9").Value = 0.001517336653189343
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.08503100000000001
$ws.Range("N9").Value = 0.255093
$ws.Range("O9").Value = 0.6282924065909707
$ws.Range("P9").Value = 0.6282924065909707
$ws.Range("Q9").Value = 0.03415343808533334
$ws.Range("R9").Value = 0.307380942768
$ws.Range("S9").Value = 0.0009533310974410213
$ws.Range("T9").Value = 0.0009533310974410213

